$wb = $excel.ActiveWorkbook

# --- NameData sheet: add 3 new rows (8,9,10) ---
$wsName = $wb.Worksheets.Item("NameData")
$wsName.Range("A8").Value = "7"
$wsName.Range("F8").Value = "Wills Company"
$wsName.Range("A9").Value = "8"
$wsName.Range("D9").Value = "Carlos Jacinta"
$wsName.Range("A10").Value = "9"
$wsName.Range("B10").Value = "Bridges"
$wsName.Range("C10").Value = "Waters"
[void]$wsName.Range("C10").Select()

# --- ACHData sheet: add 2 new rows (9,10) ---
$wsAch = $wb.Worksheets.Item("ACHData")
$wsAch.Range("A9").Value = "8"
$wsAch.Range("B9").Value = "95125480"
$wsAch.Range("C9").Value = "95125480"
$wsAch.Range("D9").Value = "256072691"
$wsAch.Range("F9").Value = "999999999"
$wsAch.Range("G9").Value = "Corporate Check"
$wsAch.Range("A10").Value = "9"
$wsAch.Range("B10").Value = "95125480"
$wsAch.Range("C10").Value = "95125480"
$wsAch.Range("D10").Value = "256072691"
$wsAch.Range("E10").Value = "1"
$wsAch.Range("G10").Value = "Personal Checking"
[void]$wsAch.Range("A10").Select()

# --- Make NameData the active sheet/tab (was UDFData) ---
[void]$wsName.Activate()
